$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# Capture existing row-2 values for G2 ("Acceptance Rate...") and the
# trailing column ("Majors Available...") before shifting them rightward.
$oldF2 = $ws.Range("F2").Value2
$oldG2 = $ws.Range("G2").Value2

# Shift G2 -> H2 (value + format) to make room for the new column.
$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial($xlPasteFormats)
$ws.Range("H2").Value = $oldG2

# Shift F2 -> G2 (value + format).
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial($xlPasteFormats)
$ws.Range("G2").Value = $oldF2

# New F2 cell with the added shared string, styled like the rest of the
# "Applying" category row (same fill as E2).
$ws.Range("E2").Copy()
$ws.Range("F2").PasteSpecial($xlPasteFormats)
$ws.Range("F2").Value = "SAT/ACT scores must be received by"

$excel.CutCopyMode = $false

# Update selection to reflect the new active cell recorded in the diff.
$ws.Range("G1").Select()
